# Fixed milestones & timeline dates
# The milestones table has "Briefing #3" moving later in the schedule
# (after "Rank Fuzzing Tools Based on Probability of Success"), so the
# rows between "Briefing #2" and "Rank Fuzzing Tools..." each shift up
# one slot, and the dates are adjusted to match the new schedule.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 16: "Briefing #2" - both dates 10/14 -> 10/16
$t.Cell(16, 2).Range.Text = "10/16"
$t.Cell(16, 3).Range.Text = "10/16"

# Row 17: was "Briefing #3" -> becomes "Test Fuzzing Tool #3"
$t.Cell(17, 1).Range.Text = "Test Fuzzing Tool #3"
$t.Cell(17, 2).Range.Text = "10/15"
$t.Cell(17, 3).Range.Text = "10/21"

# Row 18: was "Test Fuzzing Tool #3" -> becomes "Fuzzing Tool #3 Testing Finished"
$t.Cell(18, 1).Range.Text = "Fuzzing Tool #3 Testing Finished"
$t.Cell(18, 2).Range.Text = "10/21"

# Row 19: was "Fuzzing Tool #3 Testing Finished" -> becomes "Analyze Fuzz Testing Results"
$t.Cell(19, 1).Range.Text = "Analyze Fuzz Testing Results"

# Row 20: was "Analyze Fuzz Testing Results" -> becomes "Rank Fuzzing Tools Based on Probability of Success"
$t.Cell(20, 1).Range.Text = "Rank Fuzzing Tools Based on Probability of Success"
$t.Cell(20, 2).Range.Text = "10/22"
$t.Cell(20, 3).Range.Text = "10/22"

# Row 21: was "Rank Fuzzing Tools Based on Probability of Success" -> becomes "Briefing #3"
$t.Cell(21, 1).Range.Text = "Briefing #3"
$t.Cell(21, 2).Range.Text = "10/28"
$t.Cell(21, 3).Range.Text = "10/28"
